$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 889.82355
$ws.Range("I28").Value = 547.46155
$ws.Range("J28").Value = 2002.5
$ws.Range("K28").Value = 547.46155
$ws.Range("L28").Value = 2002.5
$ws.Range("M28").Value = -62.46154999999999
$ws.Range("N28").Value = -2972.5

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6716.3076
$ws.Range("I70").Value = 3938
$ws.Range("J70").Value = 7549.8
$ws.Range("K70").Value = 11814
$ws.Range("L70").Value = 22649.4
$ws.Range("M70").Value = -11544
$ws.Range("N70").Value = -23189.4

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6716.3076
$ws.Range("I73").Value = 3938
$ws.Range("J73").Value = 7549.8
$ws.Range("K73").Value = 11814
$ws.Range("L73").Value = 22649.4
$ws.Range("M73").Value = -10878
$ws.Range("N73").Value = -24521.4

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 3000
$ws.Range("N80").Value = -4996

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 9000
$ws.Range("N83").Value = -18984

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5584.143
$ws.Range("I98").Value = 2719.8647
$ws.Range("J98").Value = 26779.8
$ws.Range("K98").Value = 2719.8647
$ws.Range("L98").Value = 26779.8
$ws.Range("M98").Value = -1221.8647
$ws.Range("N98").Value = -29775.8

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 5584.143
$ws.Range("I122").Value = 2719.8647
$ws.Range("J122").Value = 26779.8
$ws.Range("K122").Value = 8159.5941
$ws.Range("L122").Value = 80339.39999999999
$ws.Range("M122").Value = -5709.5941
$ws.Range("N122").Value = -85239.39999999999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1815.2778
$ws.Range("J137").Value = 2253.4167
$ws.Range("L137").Value = 6760.250100000001
$ws.Range("N137").Value = -11860.2501

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2300.1804
$ws.Range("I32").Value = 1901.8727
$ws.Range("J32").Value = 5951.3335
$ws.Range("K32").Value = 1901.8727
$ws.Range("L32").Value = 5951.3335
$ws.Range("M32").Value = -1614.8727
$ws.Range("N32").Value = -6525.3335

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10139.389
$ws.Range("I61").Value = 7213.75
$ws.Range("K61").Value = 7213.75
$ws.Range("M61").Value = -7001.75

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2949.75
$ws.Range("I63").Value = 2933
$ws.Range("K63").Value = 2933
$ws.Range("M63").Value = -2247

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2949.75
$ws.Range("I66").Value = 2933
$ws.Range("K66").Value = 14665
$ws.Range("M66").Value = -11233

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2117.9
$ws.Range("I110").Value = 1921.1177
$ws.Range("J110").Value = 3233
$ws.Range("K110").Value = 1921.1177
$ws.Range("L110").Value = 3233
$ws.Range("M110").Value = 123.8823
$ws.Range("N110").Value = -7323

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2449.7646
$ws.Range("I132").Value = 2377.7334
$ws.Range("K132").Value = 7133.2002
$ws.Range("M132").Value = -4603.2002

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10139.389
$ws.Range("I136").Value = 7213.75
$ws.Range("K136").Value = 21641.25
$ws.Range("M136").Value = -19091.25

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1250487.9
$ws.Range("I22").Value = 272.14285
$ws.Range("K22").Value = 272.14285
$ws.Range("M22").Value = -99.14285000000001

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3128606
$ws.Range("I94").Value = 3448717.8
$ws.Range("K94").Value = 3448717.8
$ws.Range("M94").Value = -3448266.8

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8810.333000000001
$ws.Range("I134").Value = 9370.200000000001
$ws.Range("J134").Value = 7210.7144
$ws.Range("K134").Value = 28110.6
$ws.Range("L134").Value = 21632.1432
$ws.Range("M134").Value = -25575.6
$ws.Range("N134").Value = -26702.1432

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4848.091
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4848.091
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13026.75
$ws.Range("I62").Value = 6372.143
$ws.Range("J62").Value = 18202.555
$ws.Range("K62").Value = 6372.143
$ws.Range("L62").Value = 18202.555
$ws.Range("M62").Value = -5748.143
$ws.Range("N62").Value = -19450.555

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 13026.75
$ws.Range("I65").Value = 6372.143
$ws.Range("J65").Value = 18202.555
$ws.Range("K65").Value = 31860.715
$ws.Range("L65").Value = 91012.77499999999
$ws.Range("M65").Value = -28740.715
$ws.Range("N65").Value = -97252.77499999999

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3514.318
$ws.Range("I107").Value = 4159.3125
$ws.Range("K107").Value = 4159.3125
$ws.Range("M107").Value = -2239.3125

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6230.5356
$ws.Range("I134").Value = 5632.864
$ws.Range("K134").Value = 16898.592
$ws.Range("M134").Value = -14363.592

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3920
$ws.Range("I118").Value = 3920
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 11760
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -10517

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3708.2354
$ws.Range("J121").Value = 3810.0715
$ws.Range("L121").Value = 11430.2145
$ws.Range("N121").Value = -14050.2145

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 22733754
$ws.Range("J140").Value = 20000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3157
$ws.Range("I80").Value = 2985.75
$ws.Range("K80").Value = 2985.75
$ws.Range("M80").Value = -1987.75

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3157
$ws.Range("I83").Value = 2985.75
$ws.Range("K83").Value = 14928.75
$ws.Range("M83").Value = -9936.75

# GSM row 105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2077.111
$ws.Range("J46").Value = 2384.8572
$ws.Range("L46").Value = 2384.8572
$ws.Range("N46").Value = -2760.8572

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 431.55554
$ws.Range("I55").Value = 150
$ws.Range("K55").Value = 150
$ws.Range("M55").Value = 23

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2239.5
$ws.Range("I68").Value = 2063.2666
$ws.Range("J68").Value = 2533.2222
$ws.Range("K68").Value = 2063.2666
$ws.Range("L68").Value = 2533.2222
$ws.Range("M68").Value = -1314.2666
$ws.Range("N68").Value = -4031.2222

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2239.5
$ws.Range("I71").Value = 2063.2666
$ws.Range("J71").Value = 2533.2222
$ws.Range("K71").Value = 10316.333
$ws.Range("L71").Value = 12666.111
$ws.Range("M71").Value = -6572.332999999999
$ws.Range("N71").Value = -20154.111

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2103.25
$ws.Range("I82").Value = 1652
$ws.Range("K82").Value = 1652
$ws.Range("M82").Value = -1291

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2103.25
$ws.Range("I85").Value = 1652
$ws.Range("K85").Value = 1652
$ws.Range("M85").Value = -404

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 984.5714
$ws.Range("I93").Value = 955.8570999999999
$ws.Range("J93").Value = 1013.2857
$ws.Range("K93").Value = 955.8570999999999
$ws.Range("L93").Value = 1013.2857
$ws.Range("M93").Value = 292.1429000000001
$ws.Range("N93").Value = -3509.2857

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2705.25
$ws.Range("I136").Value = 1602.3572
$ws.Range("J136").Value = 4249.3
$ws.Range("K136").Value = 4807.071599999999
$ws.Range("L136").Value = 12747.9
$ws.Range("M136").Value = -2257.071599999999
$ws.Range("N136").Value = -17847.9

# WVR row 51
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 9198.333000000001
$ws.Range("I51").Value = 9198.333000000001
$ws.Range("K51").Value = 9198.333000000001
$ws.Range("M51").Value = -8688.333000000001

# WVR row 58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 28992.5
$ws.Range("I58").Value = 28992.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 28992.5
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -28684.5

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3733.8
$ws.Range("I96").Value = 1850
$ws.Range("J96").Value = 4989.6665
$ws.Range("K96").Value = 1850
$ws.Range("L96").Value = 4989.6665
$ws.Range("M96").Value = -477
$ws.Range("N96").Value = -7735.6665

# WVR row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 22500
$ws.Range("J105").Value = 22500
$ws.Range("L105").Value = 22500
$ws.Range("N105").Value = -29488

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 857
$ws.Range("I107").Value = 1003.625
$ws.Range("J107").Value = 622.4
$ws.Range("K107").Value = 3010.875
$ws.Range("L107").Value = 1867.2
$ws.Range("M107").Value = -1090.875
$ws.Range("N107").Value = -5707.2

